$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Add the new "OVERALL" column header (copy formatting from the
#    existing header cells so it matches their style: bold, filled,
#    bordered, left aligned).
# ------------------------------------------------------------------
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "OVERALL"

# ------------------------------------------------------------------
# 2. Remove the "Wagner Junior" row (row 2). This shifts the Messi
#    row up from row 3 to row 2.
# ------------------------------------------------------------------
$ws.Rows.Item(2).Delete()

# ------------------------------------------------------------------
# 3. Refresh Messi's attributes (now on row 2) and add his OVERALL
#    rating. Values are written as real numbers.
# ------------------------------------------------------------------
$ws.Range("C2").Value = 98
$ws.Range("D2").Value = 97
$ws.Range("E2").Value = 58
$ws.Range("F2").Value = 40
$ws.Range("G2").Value = 99
$ws.Range("H2").Value = 78

# ------------------------------------------------------------------
# 4. Bold the header row (all eight headers, including the new one).
# ------------------------------------------------------------------
$ws.Range("A1:H1").Font.Bold = $true

# ------------------------------------------------------------------
# 5. Fit the new column's width to its content and reset selection.
# ------------------------------------------------------------------
[void]$ws.Columns.Item(8).AutoFit()
[void]$ws.Range("A1").Select()
$excel.CutCopyMode = $false
